$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.011.76"
$ws.Range("E2").Value = "'  +0.29%  "
$ws.Range("D3").Value = "'3.085.37"
$ws.Range("E3").Value = "'  +0.37%  "
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'571.64"
$ws.Range("E5").Value = "'  -0.97%  "
$ws.Range("D6").Value = "'176.69"
$ws.Range("E6").Value = "'  +5.34%  "
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("D8").Value = "'3.083.85"
$ws.Range("E8").Value = "'  +0.39%  "
$ws.Range("E9").Value = "'  -0.31%  "
$ws.Range("D10").Value = "'6.39"
$ws.Range("E10").Value = "'  +0.16%  "
$ws.Range("E11").Value = "'  +1.20%  "
$ws.Range("E12").Value = "'  -1.08%  "
$ws.Range("E13").Value = "'  +0.18%  "
$ws.Range("D14").Value = "'35.77"
$ws.Range("E14").Value = "'  -0.77%  "
$ws.Range("D16").Value = "'3.601.92"
$ws.Range("E16").Value = "'  +0.27%  "
$ws.Range("D17").Value = "'66.980.70"
$ws.Range("D18").Value = "'6.98"
$ws.Range("E18").Value = "'  +0.04%  "
$ws.Range("D19").Value = "'3.088.66"
$ws.Range("E19").Value = "'  +0.18%  "
$ws.Range("D20").Value = "'16.41"
$ws.Range("E20").Value = "'  +0.31%  "
$ws.Range("D21").Value = "'486.58"
$ws.Range("E21").Value = "'  +0.57%  "
$ws.Range("D22").Value = "'7.64"
$ws.Range("E22").Value = "'  -0.83%  "
$ws.Range("D23").Value = "'0.682"
$ws.Range("E23").Value = "'  -0.54%  "
$ws.Range("D24").Value = "'83.25"
$ws.Range("E24").Value = "'  +0.75%  "
$ws.Range("D25").Value = "'12.56"
$ws.Range("E25").Value = "'  -1.75%  "
$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "'  +0.93%  "
$ws.Range("D27").Value = "'10.12"
$ws.Range("E27").Value = "'  -0.56%  "
$ws.Range("E28").Value = "'  +0.00%  "
$ws.Range("D29").Value = "'7.83"
$ws.Range("E29").Value = "'  +1.82%  "
$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = "'  -1.10%  "
$ws.Range("E31").Value = "'  -1.92%  "
$ws.Range("D32").Value = "'27.91"
$ws.Range("E32").Value = "'  +0.68%  "
$ws.Range("E33").Value = "'  +0.11%  "
$ws.Range("D34").Value = "'0.0₃0941"
$ws.Range("E34").Value = "'  +4.48%  "
$ws.Range("E35").Value = "'  -0.02%  "
$ws.Range("D36").Value = "'47.16"
$ws.Range("E36").Value = "'  +1.96%  "
$ws.Range("D37").Value = "'0.941"
$ws.Range("E37").Value = "'  -1.13%  "
$ws.Range("E38").Value = "'  -2.12%  "
$ws.Range("E39").Value = "'  +3.43%  "
$ws.Range("D40").Value = "'48.99"
$ws.Range("E40").Value = "'  -0.51%  "
$ws.Range("E41").Value = "'  +1.37%  "
$ws.Range("E42").Value = "'  -0.11%  "
$ws.Range("E43").Value = "'  +8.76%  "
$ws.Range("D44").Value = "'8.19"
$ws.Range("E44").Value = "'  -1.16%  "
$ws.Range("D45").Value = "'2.786.48"
$ws.Range("E45").Value = "'  +0.79%  "
$ws.Range("D46").Value = "'366.95"
$ws.Range("E46").Value = "'  -0.88%  "
$ws.Range("D47").Value = "'0.0342"
$ws.Range("E47").Value = "'  -0.67%  "
$ws.Range("D48").Value = "'134.65"
$ws.Range("E48").Value = "'  -0.60%  "
$ws.Range("D50").Value = "'25.33"
$ws.Range("E50").Value = "'  +4.28%  "
$ws.Range("E51").Value = "'  +7.79%  "
